$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 384399.8731553834
$ws.Range("C2").Value = 324789.5172530132
$ws.Range("D2").Value = -59610.35590237018
$ws.Range("E2").Value = 0.1835353443871567
